# Auto-update draw results: append the 2025-10-28 Pick 3 draw as a new
# row (42) at the bottom of the Results sheet.
#
# Columns A (Date, "2025-10-28") and C (Phase, "251028") look like a date
# / a plain number, so a bare .Value assignment would let Excel coerce
# them into a date-serial / numeric value. Prefixing with a leading
# apostrophe forces them to be stored as literal text, matching how the
# rest of the column already stores these look-like-numbers values as
# strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 42

$ws.Cells.Item($row, 1).Value = "'2025-10-28"
$ws.Cells.Item($row, 2).Value = "Pick 3"
$ws.Cells.Item($row, 3).Value = "'251028"
$ws.Cells.Item($row, 4).Value = "8-6-0"
$ws.Cells.Item($row, 5).Value = "2025-10-28T21:40:17.822+04:00"
